# This workbook contains a single sheet of weekly price records for
# "Piña" at "Terminal Hortofrutícola Agro Chillán", ordered from most
# recent date (row 2) to oldest (row 380). The edit adds one new, more
# recent weekly record right after the two newest rows (313-314), which
# pushes every existing record from row 315 onward down by one row; the
# record that used to be the very last one (row 380) is duplicated to
# become the new last row (381).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Duplicate the last existing row (380) into the new row (381).
#    Copying (rather than just setting .Value) also preserves the date
#    cell's number format/style and extends the sheet's used range /
#    dimension to A1:T381 automatically.
$ws.Range("A380:T380").Copy($ws.Range("A381:T381"))

# 2) Shift the variable columns (D, L, M, N, O, P, Q, S, T) down by one
#    row for every row from the old last row (380) up to row 316, i.e.
#    row r takes on what used to be in row r-1. Columns A, B, C, E, F,
#    G, H, I, J, K, R stay constant throughout this block, so they do
#    not need to be touched.
for ($r = 380; $r -ge 316; $r--) {
    $src = $r - 1

    $ws.Range("D$r").Value = $ws.Range("D$src").Value()
    $ws.Range("L$r").Value = $ws.Range("L$src").Value()
    $ws.Range("M$r").Value = $ws.Range("M$src").Value()
    $ws.Range("N$r").Value = $ws.Range("N$src").Value()
    $ws.Range("O$r").Value = $ws.Range("O$src").Value()
    $ws.Range("P$r").Value = $ws.Range("P$src").Value()
    $ws.Range("Q$r").Value = $ws.Range("Q$src").Value()
    $ws.Range("S$r").Value = $ws.Range("S$src").Value()
    $ws.Range("T$r").Value = $ws.Range("T$src").Value()
}

# 3) Write the brand-new record into row 315 (quality/unit/kg-per-unit
#    stay "Segunda" / "$/caja 14 unidades" / 14, only the date, volume,
#    prices and $/Kg change).
$ws.Range("D315").Value = "10/10/2023"
$ws.Range("M315").Value = 80
$ws.Range("N315").Value = 24000
$ws.Range("O315").Value = 24000
$ws.Range("P315").Value = 24000
$ws.Range("S315").Value = 1714
